$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1972.1621
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 2033.4286
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 6100.2858
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -6436.2858
$ws.Range("H123").Value = 25166.666
$ws.Range("J123").Value = 25166.666
$ws.Range("L123").Value = 25166.666
$ws.Range("N123").Value = -34966.666
$ws.Range("H127").Value = 125001050
$ws.Range("I127").Value = 200000540
$ws.Range("J127").Value = 1893.3334
$ws.Range("K127").Value = 600001620
$ws.Range("L127").Value = 5680.0002
$ws.Range("M127").Value = -599996660
$ws.Range("N127").Value = -15600.0002
$ws.Range("H129").Value = 745.0417
$ws.Range("I129").Value = 445.25
$ws.Range("J129").Value = 1344.625
$ws.Range("K129").Value = 1335.75
$ws.Range("L129").Value = 4033.875
$ws.Range("M129").Value = 3664.25
$ws.Range("N129").Value = -14033.875
$ws.Range("H130").Value = 21999.8
$ws.Range("J130").Value = 21999.8
$ws.Range("L130").Value = 21999.8
$ws.Range("N130").Value = -32039.8
$ws.Range("H131").Value = 62500664
$ws.Range("I131").Value = 62500664
$ws.Range("K131").Value = 187501992
$ws.Range("M131").Value = -187496952
$ws.Range("H132").Value = 3662.1042
$ws.Range("I132").Value = 1168.7858
$ws.Range("J132").Value = 21115.334
$ws.Range("K132").Value = 3506.3574
$ws.Range("L132").Value = 63346.00199999999
$ws.Range("M132").Value = -976.3574000000003
$ws.Range("N132").Value = -68406.00199999999
$ws.Range("H133").Value = 45353.223
$ws.Range("J133").Value = 45353.223
$ws.Range("L133").Value = 45353.223
$ws.Range("N133").Value = -55473.223
$ws.Range("H138").Value = 3335512.5
$ws.Range("I138").Value = 1926.091
$ws.Range("J138").Value = 5265483.5
$ws.Range("K138").Value = 5778.272999999999
$ws.Range("L138").Value = 15796450.5
$ws.Range("M138").Value = -638.2729999999992
$ws.Range("N138").Value = -15806730.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10102185
$ws.Range("I45").Value = 12988153
$ws.Range("K45").Value = 12988153
$ws.Range("M45").Value = -12987776
$ws.Range("H61").Value = 1464.85
$ws.Range("I61").Value = 954.2727
$ws.Range("J61").Value = 2088.889
$ws.Range("K61").Value = 954.2727
$ws.Range("L61").Value = 2088.889
$ws.Range("M61").Value = -742.2727
$ws.Range("N61").Value = -2512.889
$ws.Range("H92").Value = 21779.375
$ws.Range("J92").Value = 21779.375
$ws.Range("L92").Value = 21779.375
$ws.Range("N92").Value = -26771.375
$ws.Range("H102").Value = 2381.4285
$ws.Range("I102").Value = 1705
$ws.Range("K102").Value = 1705
$ws.Range("M102").Value = -83
$ws.Range("H136").Value = 1464.85
$ws.Range("I136").Value = 954.2727
$ws.Range("J136").Value = 2088.889
$ws.Range("K136").Value = 2862.8181
$ws.Range("L136").Value = 6266.667
$ws.Range("M136").Value = -312.8181
$ws.Range("N136").Value = -11366.667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2425.1086
$ws.Range("I105").Value = 2423
$ws.Range("J105").Value = 2431.818
$ws.Range("K105").Value = 2423
$ws.Range("L105").Value = 2431.818
$ws.Range("M105").Value = -676
$ws.Range("N105").Value = -5925.818

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2265.0652
$ws.Range("I31").Value = 1191.3334
$ws.Range("J31").Value = 2644.0293
$ws.Range("K31").Value = 1191.3334
$ws.Range("L31").Value = 2644.0293
$ws.Range("M31").Value = -896.3334
$ws.Range("N31").Value = -3234.0293
$ws.Range("H34").Value = 2265.0652
$ws.Range("I34").Value = 1191.3334
$ws.Range("J34").Value = 2644.0293
$ws.Range("K34").Value = 1191.3334
$ws.Range("L34").Value = 2644.0293
$ws.Range("M34").Value = -989.3334
$ws.Range("N34").Value = -3048.0293
$ws.Range("H99").Value = 1834.1389
$ws.Range("I99").Value = 1696.5416
$ws.Range("J99").Value = 2109.3333
$ws.Range("K99").Value = 1696.5416
$ws.Range("L99").Value = 2109.3333
$ws.Range("M99").Value = -198.5416
$ws.Range("N99").Value = -5105.3333
$ws.Range("H107").Value = 2473.875
$ws.Range("I107").Value = 998.2
$ws.Range("J107").Value = 4933.3335
$ws.Range("K107").Value = 998.2
$ws.Range("L107").Value = 4933.3335
$ws.Range("M107").Value = 921.8
$ws.Range("N107").Value = -8773.333500000001
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 32000
$ws.Range("J119").Value = 32000
$ws.Range("L119").Value = 32000
$ws.Range("N119").Value = -41676
$ws.Range("H126").Value = 1834.1389
$ws.Range("I126").Value = 1696.5416
$ws.Range("J126").Value = 2109.3333
$ws.Range("K126").Value = 5089.6248
$ws.Range("L126").Value = 6327.999899999999
$ws.Range("M126").Value = -2619.6248
$ws.Range("N126").Value = -11267.9999
$ws.Range("H141").Value = 51722.785
$ws.Range("J141").Value = 55316.848
$ws.Range("L141").Value = 55316.848
$ws.Range("N141").Value = -65676.848

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 188.84616
$ws.Range("I38").Value = 47.166668
$ws.Range("J38").Value = 310.2857
$ws.Range("K38").Value = 141.500004
$ws.Range("L38").Value = 930.8571000000001
$ws.Range("M38").Value = 205.499996
$ws.Range("N38").Value = -1624.8571
$ws.Range("H113").Value = 3367791
$ws.Range("J113").Value = 648.5
$ws.Range("L113").Value = 1945.5
$ws.Range("N113").Value = -6285.5
$ws.Range("H131").Value = 881.04
$ws.Range("I131").Value = 444.83334
$ws.Range("J131").Value = 940.5227
$ws.Range("K131").Value = 1334.50002
$ws.Range("L131").Value = 2821.5681
$ws.Range("M131").Value = 3705.49998
$ws.Range("N131").Value = -12901.5681

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 52638100
$ws.Range("I126").Value = 90915150
$ws.Range("J126").Value = 7150
$ws.Range("K126").Value = 272745450
$ws.Range("L126").Value = 21450
$ws.Range("M126").Value = -272742980
$ws.Range("N126").Value = -26390

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1062.6111
$ws.Range("I82").Value = 954.5
$ws.Range("J82").Value = 1116.6666
$ws.Range("K82").Value = 954.5
$ws.Range("L82").Value = 1116.6666
$ws.Range("M82").Value = -593.5
$ws.Range("N82").Value = -1838.6666
$ws.Range("H85").Value = 1062.6111
$ws.Range("I85").Value = 954.5
$ws.Range("J85").Value = 1116.6666
$ws.Range("K85").Value = 954.5
$ws.Range("L85").Value = 1116.6666
$ws.Range("M85").Value = 293.5
$ws.Range("N85").Value = -3612.6666
$ws.Range("H100").Value = 1850.875
$ws.Range("I100").Value = 1834.3334
$ws.Range("J100").Value = 1860.8
$ws.Range("K100").Value = 1834.3334
$ws.Range("L100").Value = 1860.8
$ws.Range("M100").Value = -1293.3334
$ws.Range("N100").Value = -2942.8
$ws.Range("H122").Value = 3837.8235
$ws.Range("I122").Value = 5040.5
$ws.Range("J122").Value = 3181.818
$ws.Range("K122").Value = 15121.5
$ws.Range("L122").Value = 9545.454000000002
$ws.Range("M122").Value = -12671.5
$ws.Range("N122").Value = -14445.454
$ws.Range("H136").Value = 6097.087
$ws.Range("I136").Value = 1468.5
$ws.Range("J136").Value = 22760
$ws.Range("K136").Value = 4405.5
$ws.Range("L136").Value = 68280
$ws.Range("M136").Value = -1855.5
$ws.Range("N136").Value = -73380
